$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 ("naman roy" row), shifting it down to row 4.
$ws.Rows.Item(3).Insert()

# Row 3: new candidate "tishya" (rejected at L1)
$ws.Cells.Item(3, 2).Value = 308
$ws.Cells.Item(3, 3).Value = "tishya"
$ws.Cells.Item(3, 4).Value = "tishya@gmail.com"
$ws.Cells.Item(3, 5).Value = "globalTiger"
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = "55667788"
$ws.Cells.Item(3, 7).Value = "L1_rejected"

# Row 5: new candidate "peter griffin" (pending - TBS)
$ws.Cells.Item(5, 2).Value = 315
$ws.Cells.Item(5, 3).Value = "peter griffin"
$ws.Cells.Item(5, 4).Value = "peter@gmail.com"
$ws.Cells.Item(5, 5).Value = "nasa"
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = "55667799"
$ws.Cells.Item(5, 7).Value = "L1_TBS"

# Row 6: new candidate "surma" (pending - scheduled)
$ws.Cells.Item(6, 2).Value = 316
$ws.Cells.Item(6, 3).Value = "surma"
$ws.Cells.Item(6, 4).Value = "surma@gmail.com"
$ws.Cells.Item(6, 5).Value = "umbrala corporation"
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = "55667788"
$ws.Cells.Item(6, 7).Value = "L1_scheduled"
